$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "value" column (B) entirely and rewrite column A contents
$ws.Range("B1:B4").Delete()

$ws.Range("A1").Value = "datasetId"
$ws.Range("A2").Value = 60
$ws.Range("A3").Value = 64
$ws.Range("A4").Value = 65

# Move the active selection to A5, matching the saved selection in the file
$ws.Range("A5").Select()
